$wb = $excel.ActiveWorkbook

# The "AddCustomerTest" sheet (second sheet in the workbook) has its
# runmode flag in A2 flipped from "Y" to "N", and the active selection
# moves from A4 to A3.
$ws = $wb.Worksheets.Item("AddCustomerTest")

$ws.Range("A2").Value = "N"

$ws.Range("A3").Select()
